$wb = $excel.ActiveWorkbook

# Add the new "configuration" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "configuration"

$newSheet.Range("A1").Value = "CT Version"
$newSheet.Range("B1").Value = "SNOMED=January 31, 2018"
$newSheet.Range("A2").Value = "CT Version"
$newSheet.Range("B2").Value = "SPONSOR =   12"

$newSheet.Range("A1").Font.Bold = $true
$newSheet.Range("A2").Font.Bold = $true

$newSheet.Columns.Item(2).ColumnWidth = 24.833333333333332

$newSheet.Activate()
[void]$newSheet.Range("B13").Select()
$excel.ActiveWindow.Zoom = 170

